$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Kampagnendaten Sternbild Orion 2022: 16.-25. Januar, 14.-23. Februar, 14.-24. März",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Kampagnendaten 2022 für das Sternbild Sternbild Orion: 16.-25. Januar, 14.-23. Februar, 14.-24. März",
    2
)
